$d = $word.ActiveDocument

# --- Locate the two target paragraphs ---------------------------------
# The template has TWO "{STRUCTURE_ADRESSE}" / "{STRUCTURE_CODE_POSTAL}"
# blocks: one in the letterhead (top of the letter) and one in the body,
# right after "... a l'adresse suivante :". The diff only touches the
# second (body) occurrence, so we anchor on that sentence and take the
# two paragraphs that follow it.

$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*adresse suivante*") {
        $anchorIndex = $i
        break
    }
}
if ($anchorIndex -eq -1) {
    throw "Could not find the 'adresse suivante' anchor paragraph"
}

$adresseIndex = -1
$codePostalIndex = -1
for ($i = $anchorIndex + 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*{STRUCTURE_ADRESSE}*") {
        $adresseIndex = $i
        $codePostalIndex = $i + 1
        break
    }
}
if ($adresseIndex -eq -1) {
    throw "Could not find the '{STRUCTURE_ADRESSE}' paragraph"
}

# --- Paragraph 1: "{STRUCTURE_ADRESSE}" -> "{STRUCTURE_COURRIER_ADRESSE}"
# Insert a new run "COURRIER_" between the existing "STRUCTURE_" run and
# the existing "ADRESSE" run, leaving every other run untouched.
$p1Xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="674E9F89" w14:textId="719CB023" w:rsidR="00252F4D" w:rsidRPr="0069142E" w:rsidRDefault="00453D4E" w:rsidP="0069142E"><w:r w:rsidRPr="00255663"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>{</w:t></w:r><w:r w:rsidR="008327AE" w:rsidRPr="00255663"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>STRUCTURE_</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>COURRIER_</w:t></w:r><w:r w:rsidRPr="00255663"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>ADRESSE</w:t></w:r><w:r w:rsidR="003A4AE6" w:rsidRPr="00255663"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$d.Paragraphs($adresseIndex).Range.InsertXML($p1Xml)

# --- Paragraph 2: "{STRUCTURE_CODE_POSTAL}, {STRUCTURE_VILLE}" ---------
# -> "{STRUCTURE_COURRIER_CODE_POSTAL}, {STRUCTURE_COURRIER_VILLE}"
# Split the "{STRUCTURE_CODE_POSTAL}" run into three runs and insert a
# new "COURRIER_" run between the second "STRUCTURE_" run and "VILLE}".
$p2Xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="0D8C77F5" w14:textId="6C0C6E6B" w:rsidR="00453D4E" w:rsidRPr="00255663" w:rsidRDefault="00616916" w:rsidP="00616916"><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr><w:r w:rsidRPr="00255663"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>{STRUCTURE_</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>COURRIER_</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>CODE_POSTAL}</w:t></w:r><w:r w:rsidR="009C7FBD"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00453D4E" w:rsidRPr="00255663"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>{</w:t></w:r><w:r w:rsidR="00255663" w:rsidRPr="00255663"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>STRUCTURE_</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>COURRIER_</w:t></w:r><w:r w:rsidR="00453D4E" w:rsidRPr="00255663"><w:rPr><w:b/><w:bCs/></w:rPr><w:t>VILLE}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$d.Paragraphs($codePostalIndex).Range.InsertXML($p2Xml)

Write-Host "Done: updated paragraphs $adresseIndex and $codePostalIndex"
